# Update countries & provincias Spain
# - Update case figures for a handful of provinces (Madrid, Cataluña, Navarra,
#   La Rioja, Asturias, Caceres, Cantabria).
# - Re-sort the data block by "Casos totales" (column B) descending, since the
#   sheet is always kept sorted by that column and the updated totals change
#   the relative order.
# - Bump the "last updated" timestamp in the title cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First + last data rows (row 3 is the header row).
$firstRow = 4
$lastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1

# New values, keyed by province name, as (TotalCases, ActiveCases, Recovered, Deaths).
$updates = @{
    "Madrid"    = @(9702, 1186, 6931, 1201)
    "Cataluña"  = @(4704, 3, 4078, 191)
    "Navarra"   = @(794, 2, 652, 14)
    "La Rioja"  = @(654, 13, 536, 18)
    "Asturias"  = @(545, 12, 467, 8)
    "Caceres"   = @(243, 2, 231, 12)
    "Cantabria" = @(282, 11, 200, 5)
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($updates.ContainsKey($name)) {
        $vals = $updates[$name]
        $ws.Cells.Item($r, 2).Value = $vals[0]
        $ws.Cells.Item($r, 3).Value = $vals[1]
        $ws.Cells.Item($r, 4).Value = $vals[2]
        $ws.Cells.Item($r, 5).Value = $vals[3]
    }
}

# Re-sort the whole data block (A..E) descending by "Casos totales" (column B),
# matching the canonical ordering used by the published sheet.
$sortRange = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, 5))
$sortKey = $ws.Range($ws.Cells.Item($firstRow, 2), $ws.Cells.Item($lastRow, 2))
$sortRange.Sort($sortKey, 2)

# Bump the "last updated" timestamp shown in the title cell.
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 13:16"
